$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Snail description (row 33, column B)
$ws.Range("B33").Value = "I'm a Slowpoke. Yep, no one cares."

# Update Caterpillar description (row 34, column B)
$ws.Range("B34").Value = "They say i'm ugly. Just wait after my metamorphosis! "

# Update the view state to match the commit (scrolled/selected cell)
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B34").Select()
